$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hotel Data")

# Update the "Total Price" column (D) for the three hotel rows with computed
# guest-capacity based totals, replacing the previous "N/A" placeholders.
$ws.Range("D2").Value = "₹ 127,909"
$ws.Range("D3").Value = "₹ 1,284,055"
$ws.Range("D4").Value = "₹ 96,194"
